{"js": "// Office.js (Word JavaScript API) script.\n// Body is `async (context) => { ... }`.\n//\n// Net change (per the commit's diff):\n//   1. Remove the blank paragraph right after the \"HUBA-025: Visualizaci\u00f3n\n//      Detallada de Comodidades Disponibles\" heading (the one that used to\n//      sit before \"Como: Comodidades\").\n//   2. Remove the blank paragraph right after the \"HUBA-026: Identificaci\u00f3n\n//      Clara de Comodidades No Disponibles\" heading (same situation).\n//   3. Append a brand-new user story block (HUBA-027) at the end of the\n//      document: a blank line, the heading, then Como/Quiero/Para\n//      paragraphs, followed by a trailing blank paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nfunction findHeadingIndex(marker) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(marker) !== -1) {\n      return i;\n    }\n  }\n  return -1;\n}\n\n// 1) Delete the empty paragraph that follows the HUBA-025 heading.\nconst idx025 = findHeadingIndex(\"HUBA-025\");\nif (idx025 !== -1 && items[idx025 + 1].text === \"\") {\n  items[idx025 + 1].delete();\n}\n\n// 2) Delete the empty paragraph that follows the HUBA-026 heading.\nconst idx026 = findHeadingIndex(\"HUBA-026\");\nif (idx026 !== -1 && items[idx026 + 1].text === \"\") {\n  items[idx026 + 1].delete();\n}\n\nawait context.sync();\n\n// 3) Append the new HUBA-027 user story block at the very end of the body.\nconst finalParagraphs = body.paragraphs;\nfinalParagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = finalParagraphs.items[finalParagraphs.items.length - 1];\nconst insertionRange = lastParagraph.getRange(\"After\");\n\nconst rPr =\n  '<w:rPr><w:rFonts w:ascii=\"Arial\" w:eastAsiaTheme=\"minorHAnsi\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n  '<w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr>';\n\nfunction run(text) {\n  const preserve = /^\\s|\\s$/.test(text) ? ' xml:space=\"preserve\"' : \"\";\n  return (\n    \"<w:r>\" +\n    rPr +\n    \"<w:t\" +\n    preserve +\n    \">\" +\n    text\n      .replace(/&/g, \"&amp;\")\n      .replace(/</g, \"&lt;\")\n      .replace(/>/g, \"&gt;\") +\n    \"</w:t></w:r>\"\n  );\n}\n\nfunction emptyParagraph() {\n  return \"<w:p><w:pPr>\" + rPr + \"</w:pPr></w:p>\";\n}\n\nfunction textParagraph(text) {\n  return \"<w:p><w:pPr>\" + rPr + \"</w:pPr>\" + run(text) + \"</w:p>\";\n}\n\nconst newBlockXml =\n  emptyParagraph() +\n  textParagraph(\n    \"HUBA-027: Administraci\u00f3n Integral de Comodidades de la Propiedad\"\n  ) +\n  textParagraph(\"Como: Comodidades\") +\n  textParagraph(\n    \"Quiero: tener la capacidad de a\u00f1adir nuevas comodidades, eliminar las existentes o actualizar la informaci\u00f3n de estas,\"\n  ) +\n  textParagraph(\n    \"Para: mantener la lista de servicios siempre al d\u00eda y as\u00ed garantizar que los usuarios accedan a la informaci\u00f3n m\u00e1s precisa y relevante, mejorando constantemente la calidad de nuestro servicio.\"\n  ) +\n  emptyParagraph();\n\nconst ooxmlPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  \"</Relationships>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  newBlockXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ninsertionRange.insertOoxml(ooxmlPackage, \"End\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# Net change (per the commit's diff):\n#   1. Remove the blank paragraph right after the \"HUBA-025: Visualizaci\u00f3n\n#      Detallada de Comodidades Disponibles\" heading (the one that used to\n#      sit before \"Como: Comodidades\").\n#   2. Remove the blank paragraph right after the \"HUBA-026: Identificaci\u00f3n\n#      Clara de Comodidades No Disponibles\" heading (same situation).\n#   3. Append a brand-new user story block (HUBA-027) at the end of the\n#      document: a blank line, the heading, then Como/Quiero/Para\n#      paragraphs, followed by a trailing blank paragraph.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($marker) {\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        if ($d.Paragraphs.Item($i).Range.Text.Contains($marker)) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# 1) Delete the empty paragraph that follows the HUBA-025 heading.\n$idx025 = Find-ParagraphIndex(\"HUBA-025\")\nif ($idx025 -ge 1) {\n    $after025 = $d.Paragraphs.Item($idx025 + 1)\n    if ($after025.Range.Text -eq [char]13) {\n        $after025.Range.Delete()\n    }\n}\n\n# 2) Delete the empty paragraph that follows the HUBA-026 heading.\n$idx026 = Find-ParagraphIndex(\"HUBA-026\")\nif ($idx026 -ge 1) {\n    $after026 = $d.Paragraphs.Item($idx026 + 1)\n    if ($after026.Range.Text -eq [char]13) {\n        $after026.Range.Delete()\n    }\n}\n\n# 3) Append the new HUBA-027 user story block at the very end of the document.\n$rPr = '<w:rPr><w:rFonts w:ascii=\"Arial\" w:eastAsiaTheme=\"minorHAnsi\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr>'\n\nfunction New-EmptyParagraphXml {\n    return \"<w:p><w:pPr>$rPr</w:pPr></w:p>\"\n}\n\nfunction New-TextParagraphXml($text) {\n    $needsPreserve = ($text -match '^\\s') -or ($text -match '\\s$')\n    $space = \"\"\n    if ($needsPreserve) { $space = ' xml:space=\"preserve\"' }\n    $escaped = $text.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n    return \"<w:p><w:pPr>$rPr</w:pPr><w:r>$rPr<w:t$space>$escaped</w:t></w:r></w:p>\"\n}\n\n$newBlock = \"\"\n$newBlock += New-EmptyParagraphXml\n$newBlock += New-TextParagraphXml(\"HUBA-027: Administraci\u00f3n Integral de Comodidades de la Propiedad\")\n$newBlock += New-TextParagraphXml(\"Como: Comodidades\")\n$newBlock += New-TextParagraphXml(\"Quiero: tener la capacidad de a\u00f1adir nuevas comodidades, eliminar las existentes o actualizar la informaci\u00f3n de estas,\")\n$newBlock += New-TextParagraphXml(\"Para: mantener la lista de servicios siempre al d\u00eda y as\u00ed garantizar que los usuarios accedan a la informaci\u00f3n m\u00e1s precisa y relevante, mejorando constantemente la calidad de nuestro servicio.\")\n$newBlock += New-EmptyParagraphXml\n\n$ooxmlPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $newBlock + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertXML($ooxmlPackage)\n"}
